$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New DEC codes to insert (968 al 1026 test range per commit message)
$newCodes = @(
    "DEC_0984", "DEC_0985", "DEC_0986", "DEC_0988", "DEC_0989", "DEC_0990",
    "DEC_0993", "DEC_0997", "DEC_0998", "DEC_0999", "DEC_1000", "DEC_1002",
    "DEC_1003", "DEC_1004", "DEC_1005", "DEC_1025", "DEC_1026"
)

# Insert 22 new rows before the old row 131 (blank separator row),
# shifting the summary/footer block (old rows 132-136) down to 154-158.
$insertRange = $ws.Range("A131:A152")
$insertRange.EntireRow.Insert()

# Seed rows 131-147 with the same formatting as the existing data block
# (rows 129-130) before writing the new values, so the text-number-format
# and quote-prefix formatting on column C carry over.
$ws.Range("A129:J130").Copy() | Out-Null
$ws.Range("A131:J147").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Rows 148-153 are blank separator rows - copy the separator formatting
# that used to live on the old row 131 (now shifted to row 153).
$ws.Range("B153:C153").Copy() | Out-Null
$ws.Range("B148:C152").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill new data rows 131-147:
#   A: DEC_xxxx code, B: "18092588-0", C: "sebA$1357", D-J: "SIN_DATO"
$row = 131
foreach ($code in $newCodes) {
    $ws.Range("A$row").Value = $code
    $ws.Range("B$row").Value = "18092588-0"
    $ws.Range("C$row").Value = "sebA`$1357"
    for ($col = 4; $col -le 10; $col++) {
        $ws.Cells.Item($row, $col).Value = "SIN_DATO"
    }
    $row++
}

# Update the saved view state to match the new scroll/selection position.
$ws.Application.ActiveWindow.ScrollRow = 130
$ws.Range("C150").Select()
